$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder: merge "A" + " " + "slide" runs into a single run "A slide".
# Setting the exact same concatenated text as a no-op leaves the existing run
# split, so first set a distinct placeholder value to force a rewrite, then
# set the final text.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "__tmp__"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "A slide"

# Caption textbox: merge "Followed"/" "/"by"/" "/"a"/" "/"picture" runs into one.
$s.Shapes.Item(4).TextFrame.TextRange.Text = "__tmp__"
$s.Shapes.Item(4).TextFrame.TextRange.Text = "Followed by a picture"
